# Update the "Förändrad" (C) column date from 45208 to 45212 for rows 2-9
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2..9) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45208) {
        $cell.Value2 = 45212
    }
}

# Update hyperlink formulas in row 2 to point to renamed files
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2305/artfynd/A 30834-2023 artfynd.xlsx", "A 30834-2023")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2305/kartor/A 30834-2023 karta.png", "A 30834-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2305/klagomål/A 30834-2023 fsc-klagomål.docx", "A 30834-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2305/klagomålsmail/A 30834-2023 fsc-klagomål mail.docx", "A 30834-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2305/tillsyn/A 30834-2023 tillsynsbegäran.docx", "A 30834-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2305/ti,llsynsmail/A 30834-2023 tillsynsbegäran mail.docx", "A 30834-2023")'
